# Update the Lama1-Itgb8 NATMI TPM results: the sending/target cluster set
# shrank from {FAPs, ECs, MuSCs, Resolving-Mac} to {FAPs, MuSCs} and every
# ligand/receptor expression + specificity metric was recomputed against the
# new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One hashtable per data row (row 1 is the unchanged header row).
$rows = @(
    @{
        Row = 2
        A = "FAPs";  B = "Lama1"; C = "Itgb8"; D = "FAPs"
        E = 3;  F = 1
        G = 0.5587383333333333;  H = 1.676215
        I = 0.8486764927018626;  J = 0.8937587278261895
        K = 3;  L = 1
        M = 2.509764;             N = 7.529292
        O = 0.4449267202378082;   P = 0.545938718598321
        Q = 1.40230135442;        R = 12.62071218978
        S = 0.3775988484407659;   T = 0.4879374946054954
    },
    @{
        Row = 3
        A = "FAPs";  B = "Lama1"; C = "Itgb8"; D = "MuSCs"
        E = 3;  F = 1
        G = 0.5587383333333333;  H = 1.676215
        I = 0.8486764927018626;  J = 0.8937587278261895
        K = 2;  L = 1
        M = 3.131084;              N = 6.262168
        O = 0.5550732797621918;    P = 0.454061281401679
        Q = 1.749456655686667;     R = 10.49673993412
        S = 0.4710776442610967;    T = 0.4058212332206941
    },
    @{
        Row = 4
        A = "MuSCs"; B = "Lama1"; C = "Itgb8"; D = "FAPs"
        E = 2;  F = 1
        G = 0.09962599999999999;  H = 0.199252
        I = 0.1513235072981373;   J = 0.1062412721738106
        K = 3;  L = 1
        M = 2.509764;              N = 7.529292
        O = 0.4449267202378082;    P = 0.545938718598321
        Q = 0.250037748264;        R = 1.500226489584
        S = 0.06732787179704229;   T = 0.0580012239928256
    },
    @{
        Row = 5
        A = "MuSCs"; B = "Lama1"; C = "Itgb8"; D = "MuSCs"
        E = 2;  F = 1
        G = 0.09962599999999999;  H = 0.199252
        I = 0.1513235072981373;   J = 0.1062412721738106
        K = 2;  L = 1
        M = 3.131084;              N = 6.262168
        O = 0.5550732797621918;    P = 0.454061281401679
        Q = 0.3119373745839999;    R = 1.247749498336
        S = 0.08399563550109505;   T = 0.04824004818098497
    }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $ws.Range($col + $r.Row).Value = $r[$col]
    }
}
